$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (2..6). Columns: A=ID(n) B=SupplierID(n) C=Type D=Quantity(n)
# E=Arrival Date F=Source. C/E/F are written as text; when the text looks like
# a plain number or an ISO date, prefix it with an apostrophe so Excel keeps
# it as text instead of auto-converting to a number/date, then clear the
# resulting "quote prefix" formatting so no extra style is left behind.
$rows = @(
    @{ Row = 2; A = 1; B = 1; C = "1";    D = 10;  E = "2024-06-02"; F = "test" },
    @{ Row = 3; A = 2; B = 7; C = "2";    D = 300; E = "2024-02-06"; F = "test" },
    @{ Row = 4; A = 3; B = 1; C = "cm2";  D = 200; E = "2";          F = "2" },
    @{ Row = 5; A = 4; B = 2; C = "cm2";  D = 200; E = "2";          F = "2" },
    @{ Row = 6; A = 5; B = 1; C = "2";    D = 2;   E = "2";          F = "2" }
)

function Set-TextCell($range, [string]$text) {
    $range.Value = "'" + $text
    $range.ClearFormats()
}

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = $r.B
    Set-TextCell $ws.Range("C$n") $r.C
    $ws.Range("D$n").Value = $r.D
    Set-TextCell $ws.Range("E$n") $r.E
    Set-TextCell $ws.Range("F$n") $r.F
}
